$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 17.41485066666667
$ws.Range("H2").Value = 52.24455200000001
$ws.Range("I2").Value = 0.1047285618770465
$ws.Range("J2").Value = 0.1047285618770465
$ws.Range("M2").Value = 168.1098273333333
$ws.Range("N2").Value = 504.329482
$ws.Range("O2").Value = 0.2984182258032519
$ws.Range("P2").Value = 0.298418225803252
$ws.Range("Q2").Value = 2927.607538609119
$ws.Range("R2").Value = 26348.46784748206
$ws.Range("S2").Value = 0.03125291162627432
$ws.Range("T2").Value = 0.03125291162627432
$ws.Range("G3").Value = 17.41485066666667
$ws.Range("H3").Value = 52.24455200000001
$ws.Range("I3").Value = 0.1047285618770465
$ws.Range("J3").Value = 0.1047285618770465
$ws.Range("O3").Value = 0.2893586437755394
$ws.Range("P3").Value = 0.2893586437755394
$ws.Range("Q3").Value = 2838.729251870474
$ws.Range("R3").Value = 25548.56326683427
$ws.Range("S3").Value = 0.03030411462930485
$ws.Range("T3").Value = 0.03030411462930485
$ws.Range("G4").Value = 17.41485066666667
$ws.Range("H4").Value = 52.24455200000001
$ws.Range("I4").Value = 0.1047285618770465
$ws.Range("J4").Value = 0.1047285618770465
$ws.Range("M4").Value = 165.99353
$ws.Range("N4").Value = 497.98059
$ws.Range("O4").Value = 0.294661504941043
$ws.Range("P4").Value = 0.294661504941043
$ws.Range("Q4").Value = 2890.752536582854
$ws.Range("R4").Value = 26016.77282924568
$ws.Range("S4").Value = 0.03085947565300167
$ws.Range("T4").Value = 0.03085947565300167
$ws.Range("G5").Value = 17.41485066666667
$ws.Range("H5").Value = 52.24455200000001
$ws.Range("I5").Value = 0.1047285618770465
$ws.Range("J5").Value = 0.1047285618770465
$ws.Range("M5").Value = 66.22673433333334
$ws.Range("N5").Value = 198.680203
$ws.Range("O5").Value = 0.1175616254801657
$ws.Range("P5").Value = 0.1175616254801657
$ws.Range("Q5").Value = 1153.328688556006
$ws.Range("R5").Value = 10379.95819700406
$ws.Range("S5").Value = 0.0123120599684657
$ws.Range("T5").Value = 0.0123120599684657
$ws.Range("I6").Value = 0.1785014126970782
$ws.Range("J6").Value = 0.1785014126970782
$ws.Range("M6").Value = 168.1098273333333
$ws.Range("N6").Value = 504.329482
$ws.Range("O6").Value = 0.2984182258032519
$ws.Range("P6").Value = 0.298418225803252
$ws.Range("Q6").Value = 4989.871646264613
$ws.Range("R6").Value = 44908.84481638151
$ws.Range("S6").Value = 0.05326807488043614
$ws.Range("T6").Value = 0.05326807488043615
$ws.Range("I7").Value = 0.1785014126970782
$ws.Range("J7").Value = 0.1785014126970782
$ws.Range("O7").Value = 0.2893586437755394
$ws.Range("P7").Value = 0.2893586437755394
$ws.Range("S7").Value = 0.0516509266900444
$ws.Range("T7").Value = 0.0516509266900444
$ws.Range("I8").Value = 0.1785014126970782
$ws.Range("J8").Value = 0.1785014126970782
$ws.Range("M8").Value = 165.99353
$ws.Range("N8").Value = 497.98059
$ws.Range("O8").Value = 0.294661504941043
$ws.Range("P8").Value = 0.294661504941043
$ws.Range("Q8").Value = 4927.055258750714
$ws.Range("R8").Value = 44343.49732875642
$ws.Range("S8").Value = 0.05259749489942325
$ws.Range("T8").Value = 0.05259749489942326
$ws.Range("I9").Value = 0.1785014126970782
$ws.Range("J9").Value = 0.1785014126970782
$ws.Range("M9").Value = 66.22673433333334
$ws.Range("N9").Value = 198.680203
$ws.Range("O9").Value = 0.1175616254801657
$ws.Range("P9").Value = 0.1175616254801657
$ws.Range("Q9").Value = 1965.756012700835
$ws.Range("R9").Value = 17691.80411430751
$ws.Range("S9").Value = 0.02098491622717439
$ws.Range("T9").Value = 0.02098491622717439
$ws.Range("G10").Value = 84.03051233333333
$ws.Range("H10").Value = 252.091537
$ws.Range("I10").Value = 0.5053385113032314
$ws.Range("J10").Value = 0.5053385113032314
$ws.Range("M10").Value = 168.1098273333333
$ws.Range("N10").Value = 504.329482
$ws.Range("O10").Value = 0.2984182258032519
$ws.Range("P10").Value = 0.298418225803252
$ws.Range("Q10").Value = 14126.3549190882
$ws.Range("R10").Value = 127137.1942717938
$ws.Range("S10").Value = 0.1508022219731669
$ws.Range("T10").Value = 0.1508022219731669
$ws.Range("G11").Value = 84.03051233333333
$ws.Range("H11").Value = 252.091537
$ws.Range("I11").Value = 0.5053385113032314
$ws.Range("J11").Value = 0.5053385113032314
$ws.Range("O11").Value = 0.2893586437755394
$ws.Range("P11").Value = 0.2893586437755394
$ws.Range("Q11").Value = 13697.49749659807
$ws.Range("R11").Value = 123277.4774693826
$ws.Range("S11").Value = 0.1462240662782531
$ws.Range("T11").Value = 0.1462240662782531
$ws.Range("G12").Value = 84.03051233333333
$ws.Range("H12").Value = 252.091537
$ws.Range("I12").Value = 0.5053385113032314
$ws.Range("J12").Value = 0.5053385113032314
$ws.Range("M12").Value = 165.99353
$ws.Range("N12").Value = 497.98059
$ws.Range("O12").Value = 0.294661504941043
$ws.Range("P12").Value = 0.294661504941043
$ws.Range("Q12").Value = 13948.52136991854
$ws.Range("R12").Value = 125536.6923292668
$ws.Range("S12").Value = 0.1489038062452764
$ws.Range("T12").Value = 0.1489038062452764
$ws.Range("G13").Value = 84.03051233333333
$ws.Range("H13").Value = 252.091537
$ws.Range("I13").Value = 0.5053385113032314
$ws.Range("J13").Value = 0.5053385113032314
$ws.Range("M13").Value = 66.22673433333334
$ws.Range("N13").Value = 198.680203
$ws.Range("O13").Value = 0.1175616254801657
$ws.Range("P13").Value = 0.1175616254801657
$ws.Range("Q13").Value = 5565.066416193557
$ws.Range("R13").Value = 50085.59774574201
$ws.Range("S13").Value = 0.05940841680653495
$ws.Range("T13").Value = 0.05940841680653495
$ws.Range("G14").Value = 35.158014
$ws.Range("H14").Value = 105.474042
$ws.Range("I14").Value = 0.2114315141226439
$ws.Range("J14").Value = 0.2114315141226439
$ws.Range("M14").Value = 168.1098273333333
$ws.Range("N14").Value = 504.329482
$ws.Range("O14").Value = 0.2984182258032519
$ws.Range("P14").Value = 0.298418225803252
$ws.Range("Q14").Value = 5910.407662922916
$ws.Range("R14").Value = 53193.66896630624
$ws.Range("S14").Value = 0.06309501732337459
$ws.Range("T14").Value = 0.06309501732337461
$ws.Range("G15").Value = 35.158014
$ws.Range("H15").Value = 105.474042
$ws.Range("I15").Value = 0.2114315141226439
$ws.Range("J15").Value = 0.2114315141226439
$ws.Range("O15").Value = 0.2893586437755394
$ws.Range("P15").Value = 0.2893586437755394
$ws.Range("Q15").Value = 5730.975515655966
$ws.Range("R15").Value = 51578.7796409037
$ws.Range("S15").Value = 0.06117953617793705
$ws.Range("T15").Value = 0.06117953617793705
$ws.Range("G16").Value = 35.158014
$ws.Range("H16").Value = 105.474042
$ws.Range("I16").Value = 0.2114315141226439
$ws.Range("J16").Value = 0.2114315141226439
$ws.Range("M16").Value = 165.99353
$ws.Range("N16").Value = 497.98059
$ws.Range("O16").Value = 0.294661504941043
$ws.Range("P16").Value = 0.294661504941043
$ws.Range("Q16").Value = 5836.00285164942
$ws.Range("R16").Value = 52524.02566484478
$ws.Range("S16").Value = 0.06230072814334163
$ws.Range("T16").Value = 0.06230072814334163
$ws.Range("G17").Value = 35.158014
$ws.Range("H17").Value = 105.474042
$ws.Range("I17").Value = 0.2114315141226439
$ws.Range("J17").Value = 0.2114315141226439
$ws.Range("M17").Value = 66.22673433333334
$ws.Range("N17").Value = 198.680203
$ws.Range("O17").Value = 0.1175616254801657
$ws.Range("P17").Value = 0.1175616254801657
$ws.Range("Q17").Value = 2328.400452865614
$ws.Range("R17").Value = 20955.60407579052
$ws.Range("S17").Value = 0.02485623247799062
$ws.Range("T17").Value = 0.02485623247799062
